# Add four new response rows (11-14) to the "Admin" sheet.
# Column A: sequential id, Column B/C: left blank, Column D: "respondente".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admin")

for ($r = 11; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = $r
    # Touch B/C with a no-op format so Excel keeps them as present-but-blank cells
    # (matching the source rows' shape), rather than omitting them entirely.
    $ws.Cells.Item($r, 2).Font.Bold = $false
    $ws.Cells.Item($r, 3).Font.Bold = $false
    $ws.Cells.Item($r, 4).Value = "respondente"
}
